# Case_3_218 (380 kV) voltage-magnitude results: Vmax raised from 1.05 to 1.02 p.u.,
# which shifts the bus-voltage solution reported in columns B:N (column H stays blank,
# column G stays the fixed 1.0 slack value, column A is the row index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ "B" = 1.02; "C" = 1.02985163381283; "D" = 1.0334098029751; "E" = 1.029627399942429; "F" = 1.039458241156218; "I" = 1.035472387700199; "J" = 1.034996120801076; "K" = 1.036212233826857; "L" = 1.032440761773479; "M" = 1.04224337180847; "N" = 1.015636665658253 }
    3 = @{ "B" = 1.02; "C" = 1.030775563345503; "D" = 1.034100561552578; "E" = 1.030411892956195; "F" = 1.04165550541496; "I" = 1.035746541569345; "J" = 1.035561387523069; "K" = 1.036712282452409; "L" = 1.033033503614254; "M" = 1.044247206865055; "N" = 1.015827145815178 }
    4 = @{ "B" = 1.02; "C" = 1.031373475103234; "D" = 1.03454749355653; "E" = 1.030919934860554; "F" = 1.043072034178141; "I" = 1.035922491905623; "J" = 1.035926607291227; "K" = 1.037035137052532; "L" = 1.033416814773073; "M" = 1.045538245007781; "N" = 1.015950139453989 }
    5 = @{ "B" = 1.02; "C" = 1.031624853792447; "D" = 1.034735375537411; "E" = 1.031133616663815; "F" = 1.043666317701273; "I" = 1.035996116615918; "J" = 1.036080015454715; "K" = 1.037170695587925; "L" = 1.033577903282417; "M" = 1.046079691350235; "N" = 1.016001783844458 }
    6 = @{ "B" = 1.02; "C" = 1.031667062364672; "D" = 1.034766921269561; "E" = 1.031169500695908; "F" = 1.043766029575443; "I" = 1.036008458344899; "J" = 1.036105765734785; "K" = 1.037193446535042; "L" = 1.033604947501212; "M" = 1.046170526821487; "N" = 1.016010451522588 }
    7 = @{ "B" = 1.02; "C" = 1.031376833971935; "D" = 1.034550004076983; "E" = 1.030922789690695; "F" = 1.043079979799719; "I" = 1.03592347703596; "J" = 1.035928657649874; "K" = 1.0370369490579; "L" = 1.033418967463561; "M" = 1.045545484937002; "N" = 1.015950829772962 }
    8 = @{ "B" = 1.02; "C" = 1.030163866386867; "D" = 1.033643255407682; "E" = 1.029892435469538; "F" = 1.040201922487479; "I" = 1.035565339357897; "J" = 1.035187268547696; "K" = 1.036381374985654; "L" = 1.032641129877842; "M" = 1.042921749646494; "N" = 1.015701093137519 }
    9 = @{ "B" = 1.02; "C" = 1.028026960110312; "D" = 1.032045172340095; "E" = 1.028080060486594; "F" = 1.035088808615581; "I" = 1.034923126276459; "J" = 1.033876648684567; "K" = 1.035220705477359; "L" = 1.031268692439332; "M" = 1.038254361684598; "N" = 1.015259032983073 }
    10 = @{ "B" = 1.02; "C" = 1.026602658346745; "D" = 1.030979582447574; "E" = 1.026873992441976; "F" = 1.031650096418278; "I" = 1.034487418971757; "J" = 1.033000050048301; "K" = 1.034443214804846; "L" = 1.030352517210957; "M" = 1.035111309019249; "N" = 1.01496298142846 }
    11 = @{ "B" = 1.02; "C" = 1.025985983216603; "D" = 1.030518117782608; "E" = 1.026352268554195; "F" = 1.030153539387335; "I" = 1.034296939317539; "J" = 1.032619788982021; "K" = 1.034105663268573; "L" = 1.029955510248439; "M" = 1.033742452928808; "N" = 1.014834467295868 }
    12 = @{ "B" = 1.02; "C" = 1.025756930309732; "D" = 1.030346700134242; "E" = 1.026158553750238; "F" = 1.029596474968001; "I" = 1.034225912197662; "J" = 1.032478438946442; "K" = 1.033980146541646; "L" = 1.029807999081162; "M" = 1.033232776662302; "N" = 1.014786682893652 }
    13 = @{ "B" = 1.02; "C" = 1.025806062601753; "D" = 1.030383470223525; "E" = 1.026200102766558; "F" = 1.029716021047707; "I" = 1.034241160205713; "J" = 1.032508763734326; "K" = 1.034007076432083; "L" = 1.029839642776321; "M" = 1.033342159845636; "N" = 1.014796935009217 }
    14 = @{ "B" = 1.02; "C" = 1.025967049479997; "D" = 1.030503948535666; "E" = 1.026336254464988; "F" = 1.030107516496406; "I" = 1.034291073803441; "J" = 1.032608107073062; "K" = 1.034095290777227; "L" = 1.029943317852321; "M" = 1.033700348073374; "N" = 1.014830518413355 }
    15 = @{ "B" = 1.02; "C" = 1.026066239802198; "D" = 1.030578177974424; "E" = 1.026420152174764; "F" = 1.030348572689988; "I" = 1.034321790792544; "J" = 1.03266930195911; "K" = 1.034149624622814; "L" = 1.030007189496785; "M" = 1.033920876618235; "N" = 1.014851203820721 }
    16 = @{ "B" = 1.02; "C" = 1.026643586170997; "D" = 1.031010207087834; "E" = 1.026908628301382; "F" = 1.031749255141533; "I" = 1.034500022081056; "J" = 1.033025272161045; "K" = 1.034465598099044; "L" = 1.030378858957092; "M" = 1.035201986101602; "N" = 1.014971503692269 }
    17 = @{ "B" = 1.02; "C" = 1.027005754939567; "D" = 1.03128119185671; "E" = 1.027215173308288; "F" = 1.032625811751418; "I" = 1.034611334581524; "J" = 1.033248378113699; "K" = 1.034663560219297; "L" = 1.03061191753479; "M" = 1.036003453508546; "N" = 1.015046878332369 }
    18 = @{ "B" = 1.02; "C" = 1.027217007408971; "D" = 1.031439247206571; "E" = 1.027394025194716; "F" = 1.033136364467597; "I" = 1.034676086228346; "J" = 1.033378445706701; "K" = 1.034778942082043; "L" = 1.03074782803445; "M" = 1.036470176753287; "N" = 1.015090812043071 }
    19 = @{ "B" = 1.02; "C" = 1.027289040008914; "D" = 1.031493139081752; "E" = 1.027455017444469; "F" = 1.033310327449044; "I" = 1.034698135226566; "J" = 1.033422784165665; "K" = 1.034818269739181; "L" = 1.030794165160852; "M" = 1.036629189644239; "N" = 1.015105787042916 }
    20 = @{ "B" = 1.02; "C" = 1.026966897090094; "D" = 1.031252118326418; "E" = 1.027182278824925; "F" = 1.032531841095129; "I" = 1.034599409922586; "J" = 1.033224447823347; "K" = 1.034642329668827; "L" = 1.030586915519455; "M" = 1.035917542422199; "N" = 1.015038794558901 }
    21 = @{ "B" = 1.02; "C" = 1.025919642651511; "D" = 1.030468470931915; "E" = 1.026296159068162; "F" = 1.029992263637467; "I" = 1.034276383084501; "J" = 1.032578855817097; "K" = 1.034069317586903; "L" = 1.029912789375553; "M" = 1.033594904599834; "N" = 1.014820630279049 }
    22 = @{ "B" = 1.02; "C" = 1.025261236286674; "D" = 1.029975707627982; "E" = 1.025739463931873; "F" = 1.028388703855767; "I" = 1.03407169425912; "J" = 1.032172343548074; "K" = 1.033708260279834; "L" = 1.029488678294044; "M" = 1.032127480534293; "N" = 1.014683180855268 }
    23 = @{ "B" = 1.02; "C" = 1.025610266075621; "D" = 1.030236935996525; "E" = 1.026034536581324; "F" = 1.029239441483336; "I" = 1.034180354888869; "J" = 1.032387900810639; "K" = 1.033899738011237; "L" = 1.029713532529651; "M" = 1.032906074060473; "N" = 1.014756072081268 }
    24 = @{ "B" = 1.02; "C" = 1.026984455262954; "D" = 1.031265255420596; "E" = 1.027197142273898; "F" = 1.032574304637558; "I" = 1.034604798703656; "J" = 1.033235261096748; "K" = 1.034651923106996; "L" = 1.03059821294275; "M" = 1.035956364288432; "N" = 1.015042447364193 }
    25 = @{ "B" = 1.02; "C" = 1.02857934584802; "D" = 1.032458349320216; "E" = 1.028548218401824; "F" = 1.036415805361666; "I" = 1.035090480985025; "J" = 1.034215975743081; "K" = 1.035521417627988; "L" = 1.031623713744943; "M" = 1.039466400600853; "N" = 1.015373552735334 }
}

foreach ($rowKey in $newValues.Keys) {
    $rowData = $newValues[$rowKey]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowKey").Value = $rowData[$col]
    }
}
